$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Step 2: Log in as the appropriate role"
$ws.Range("D3").Value = "I am redirected to the user's dashboard"

$ws.Range("C4").Value = 'Step 3:  Go to the "KPIs" page'
$ws.Range("D4").Value = "I am redirected to the KPI's page where a list of kpis are displayed"

$ws.Range("C5").Value = "Step 4: delete one that belongs to someone in my district"
$ws.Range("D5").Value = "The data is removed from the database."

$ws.Range("C6").Value = "Step 5: While logged in try to delete a kpi about me"
$ws.Range("D6").Value = "I am denied access to this"

$ws.Range("C7").Value = "Step 6: Try to delete a kpi of someone who is not in my district"
$ws.Range("D7").Value = "I am denied access to this"

try {
    $excel.ActiveWindow.ScrollRow = 2
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}

$ws.Range("E4").Select()
